$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("v1.0")

# Set the "connect" (S) column to 1 for the branch instruction rows 17-27.
# Downstream X (opcode value) and Y (hex) columns are formula-driven and
# will recalc automatically.
$ws.Range("S17:S27").Value = 1

# Update the view state: frozen-pane scroll position and active selection.
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("S32").Select()

$wb.Windows.Item(1).WindowState = -4143
